$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 347.4
$ws.Range("I8").Value = 15.071428
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 45.214284
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 93.78571600000001
$ws.Range("N8").Value = -15278
$ws.Range("H40").Value = 2137.25
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 2474.5
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 2474.5
$ws.Range("M40").Value = -1625
$ws.Range("N40").Value = -2824.5
$ws.Range("H112").Value = 1097.1632
$ws.Range("I112").Value = 825.6667
$ws.Range("J112").Value = 1114.8695
$ws.Range("K112").Value = 2477.0001
$ws.Range("L112").Value = 3344.6085
$ws.Range("M112").Value = -1369.0001
$ws.Range("N112").Value = -5560.6085
$ws.Range("H116").Value = 3634.75
$ws.Range("I116").Value = 4850.8335
$ws.Range("J116").Value = 2418.6667
$ws.Range("K116").Value = 4850.8335
$ws.Range("L116").Value = 2418.6667
$ws.Range("M116").Value = -1408.8335
$ws.Range("N116").Value = -9302.6667
$ws.Range("H129").Value = 822.14636
$ws.Range("I129").Value = 490
$ws.Range("J129").Value = 879.0857
$ws.Range("K129").Value = 1470
$ws.Range("L129").Value = 2637.2571
$ws.Range("M129").Value = 3530
$ws.Range("N129").Value = -12637.2571
$ws.Range("H137").Value = 793
$ws.Range("I137").Value = 788.375
$ws.Range("J137").Value = 830
$ws.Range("K137").Value = 2365.125
$ws.Range("L137").Value = 2490
$ws.Range("M137").Value = 184.875
$ws.Range("N137").Value = -7590
$ws.Range("H138").Value = 3490.7144
$ws.Range("I138").Value = 2141.3
$ws.Range("K138").Value = 6423.900000000001
$ws.Range("M138").Value = -1283.900000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 4000
$ws.Range("I11").Value = 4000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -3856
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value = 6893.1323
$ws.Range("I32").Value = 5633.873
$ws.Range("J32").Value = 22759.8
$ws.Range("K32").Value = 5633.873
$ws.Range("L32").Value = 22759.8
$ws.Range("M32").Value = -5346.873
$ws.Range("N32").Value = -23333.8
$ws.Range("H74").Value = 7297.5713
$ws.Range("I74").Value = 12493.5
$ws.Range("J74").Value = 2574
$ws.Range("K74").Value = 12493.5
$ws.Range("L74").Value = 2574
$ws.Range("M74").Value = -11619.5
$ws.Range("N74").Value = -4322
$ws.Range("H77").Value = 7297.5713
$ws.Range("I77").Value = 12493.5
$ws.Range("J77").Value = 2574
$ws.Range("K77").Value = 62467.5
$ws.Range("L77").Value = 12870
$ws.Range("M77").Value = -58099.5
$ws.Range("N77").Value = -21606
$ws.Range("H97").Value = 709.2857
$ws.Range("I97").Value = 618.75
$ws.Range("J97").Value = 830
$ws.Range("K97").Value = 618.75
$ws.Range("L97").Value = 830
$ws.Range("M97").Value = -122.75
$ws.Range("N97").Value = -1822
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("N105").Value = -16988
$ws.Range("H132").Value = 2197.56
$ws.Range("I132").Value = 1565.6666
$ws.Range("J132").Value = 2780.8462
$ws.Range("K132").Value = 4696.9998
$ws.Range("L132").Value = 8342.5386
$ws.Range("M132").Value = -2166.9998
$ws.Range("N132").Value = -13402.5386
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1780.862
$ws.Range("I134").Value = 1626.9166
$ws.Range("J134").Value = 2519.8
$ws.Range("K134").Value = 4880.7498
$ws.Range("L134").Value = 7559.400000000001
$ws.Range("M134").Value = -2345.7498
$ws.Range("N134").Value = -12629.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2727.25
$ws.Range("J3").Value = 2727.25
$ws.Range("L3").Value = 2727.25
$ws.Range("N3").Value = -2953.25
$ws.Range("H31").Value = 1761.8119
$ws.Range("I31").Value = 971.5714
$ws.Range("J31").Value = 3546.2258
$ws.Range("K31").Value = 971.5714
$ws.Range("L31").Value = 3546.2258
$ws.Range("M31").Value = -676.5714
$ws.Range("N31").Value = -4136.2258
$ws.Range("H34").Value = 1761.8119
$ws.Range("I34").Value = 971.5714
$ws.Range("J34").Value = 3546.2258
$ws.Range("K34").Value = 971.5714
$ws.Range("L34").Value = 3546.2258
$ws.Range("M34").Value = -769.5714
$ws.Range("N34").Value = -3950.2258
$ws.Range("H122").Value = 3426.923
$ws.Range("I122").Value = 2359.5557
$ws.Range("J122").Value = 5828.5
$ws.Range("K122").Value = 7078.6671
$ws.Range("L122").Value = 17485.5
$ws.Range("M122").Value = -4628.6671
$ws.Range("N122").Value = -22385.5
$ws.Range("H132").Value = 2539.16
$ws.Range("I132").Value = 1680.6666
$ws.Range("K132").Value = 5041.9998
$ws.Range("M132").Value = -2511.9998
$ws.Range("H138").Value = 34834.5
$ws.Range("J138").Value = 34834.5
$ws.Range("L138").Value = 34834.5
$ws.Range("N138").Value = -45114.5
$ws.Range("H140").Value = 68034.75
$ws.Range("J140").Value = 68034.75
$ws.Range("L140").Value = 68034.75
$ws.Range("N140").Value = -78394.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 598.25
$ws.Range("I60").Value = 520
$ws.Range("K60").Value = 1560
$ws.Range("M60").Value = -1309
$ws.Range("H109").Value = 5664.75
$ws.Range("I109").Value = 5274
$ws.Range("K109").Value = 15822
$ws.Range("M109").Value = -14782
$ws.Range("H114").Value = 284.2414
$ws.Range("I114").Value = 263.3
$ws.Range("J114").Value = 295.26315
$ws.Range("K114").Value = 789.9000000000001
$ws.Range("L114").Value = 885.78945
$ws.Range("M114").Value = 2464.1
$ws.Range("N114").Value = -7393.78945
$ws.Range("H131").Value = 854.42
$ws.Range("J131").Value = 858
$ws.Range("L131").Value = 2574
$ws.Range("N131").Value = -12654
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4444889
$ws.Range("I7").Value = 5714286
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 5714286
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -5714174
$ws.Range("N7").Value = -2224
$ws.Range("H8").Value = 4444889
$ws.Range("I8").Value = 5714286
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 5714286
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = -5714147
$ws.Range("N8").Value = -2278
$ws.Range("H102").Value = 3233.647
$ws.Range("I102").Value = 3295.8
$ws.Range("J102").Value = 3144.8572
$ws.Range("K102").Value = 3295.8
$ws.Range("L102").Value = 3144.8572
$ws.Range("M102").Value = -1673.8
$ws.Range("N102").Value = -6388.8572
$ws.Range("H126").Value = 2154.516
$ws.Range("I126").Value = 1912.1177
$ws.Range("J126").Value = 2448.8572
$ws.Range("K126").Value = 5736.3531
$ws.Range("L126").Value = 7346.571599999999
$ws.Range("M126").Value = -3266.3531
$ws.Range("N126").Value = -12286.5716
$ws.Range("H135").Value = 39006.445
$ws.Range("J135").Value = 39006.445
$ws.Range("L135").Value = 39006.445
$ws.Range("N135").Value = -49146.445
$ws.Range("H140").Value = 39039876
$ws.Range("J140").Value = 39039876
$ws.Range("L140").Value = 39039876
$ws.Range("N140").Value = -39050236
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 29004
$ws.Range("J12").Value = 29004
$ws.Range("L12").Value = 29004
$ws.Range("N12").Value = -29344
$ws.Range("H46").Value = 60283.117
$ws.Range("J46").Value = 1306.7778
$ws.Range("L46").Value = 1306.7778
$ws.Range("N46").Value = -1682.7778
$ws.Range("H53").Value = 7999.6665
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 7999.6665
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 7999.6665
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -9035.666499999999
$ws.Range("H55").Value = 308.75
$ws.Range("I55").Value = 171.73334
$ws.Range("J55").Value = 466.84616
$ws.Range("K55").Value = 171.73334
$ws.Range("L55").Value = 466.84616
$ws.Range("M55").Value = 1.266660000000002
$ws.Range("N55").Value = -812.8461600000001
$ws.Range("H132").Value = 9295.821
$ws.Range("I132").Value = 7589.8604
$ws.Range("J132").Value = 14938.615
$ws.Range("K132").Value = 22769.5812
$ws.Range("L132").Value = 44815.845
$ws.Range("M132").Value = -20239.5812
$ws.Range("N132").Value = -49875.845
$ws.Range("H133").Value = 32179.75
$ws.Range("J133").Value = 32179.75
$ws.Range("L133").Value = 32179.75
$ws.Range("N133").Value = -37239.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 32200
$ws.Range("J108").Value = 32200
$ws.Range("L108").Value = 32200
$ws.Range("N108").Value = -39880

Write-Output "done"